$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 3
$ws.Range("F5").Value = 3
$ws.Range("H5").Value = 46

$ws.Range("D5").Select()
